$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "MY BETS" row (row 7), with a new date/time label and probables (rows 8-10)
$ws.Range("A7").Value = "04 tue jun 2019 0"
$ws.Range("B7").Value = 240
$ws.Range("C7").Value = 5

$ws.Range("B8").Value = 244
$ws.Range("C8").Value = 5

$ws.Range("B9").Value = 246
$ws.Range("C9").Value = 5

$ws.Range("B10").Value = 224
$ws.Range("C10").Value = 5

$ws.Range("C11").Select()
